$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text "Code" -> "Book Code"
$ws.Range("D1").Value = "Book Code"

# Row 3 (book row 2): No column 4 -> 3
$ws.Range("E3").Value = 3

# Row 4 (book row 3): Name/Author "hel" -> "hello", Code "hi]" -> "hi", No 2 -> 3
$ws.Range("B4").Value = "hello"
$ws.Range("C4").Value = "hello"
$ws.Range("D4").Value = "hi"
$ws.Range("E4").Value = 3
